$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Append a new log row (row 39), mirroring the prior row's values/shape
# with an updated "Hora Consulta" timestamp.
$ws.Cells.Item(39, 1).Value = "30/12"
$ws.Cells.Item(39, 2).Value = "Mercado cerrado"
$ws.Cells.Item(39, 3).Value = "01/01/2023 03:40"
$ws.Cells.Item(39, 4).Value = "848,25"

# Match row 38's (unstyled) formatting rather than inheriting the
# column defaults Excel would otherwise apply to a freshly-used row.
$ws.Range("A38:D38").Copy()
$ws.Range("A39:D39").PasteSpecial(-4122)
